$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Me2)")

# "Second Me did not work out" -> rename the sheet (drop the stray close-paren)
$ws.Name = "Me2"

# The old HSTACK/VSTACK spill (rows 16:25) is being replaced by a smaller
# C3:C7 & D2:I2 spill (rows 13:17). Clear the old array's footprint first
# (you cannot touch part of a spilled array) and drop the now-unused rows.
$ws.Range("D16:G25").ClearContents()
$ws.Rows("16:25").EntireRow.Delete()

# Seed the spill body with the literal results so the sheet matches even if
# this host can't broadcast "Range & Range" the way Excel's dynamic arrays do.
$ws.Range("D13").Value = "APayment"
$ws.Range("E13").Value = "ADiscount"
$ws.Range("F13").Value = "AQuantity"
$ws.Range("G13").Value = "APayment"
$ws.Range("H13").Value = "ADiscount"

$ws.Range("C14").Value = "BQuantity"
$ws.Range("D14").Value = "BPayment"
$ws.Range("E14").Value = "BDiscount"
$ws.Range("F14").Value = "BQuantity"
$ws.Range("G14").Value = "BPayment"
$ws.Range("H14").Value = "BDiscount"

$ws.Range("C15").Value = "CQuantity"
$ws.Range("D15").Value = "CPayment"
$ws.Range("E15").Value = "CDiscount"
$ws.Range("F15").Value = "CQuantity"
$ws.Range("G15").Value = "CPayment"
$ws.Range("H15").Value = "CDiscount"

$ws.Range("C16").Value = "DQuantity"
$ws.Range("D16").Value = "DPayment"
$ws.Range("E16").Value = "DDiscount"
$ws.Range("F16").Value = "DQuantity"
$ws.Range("G16").Value = "DPayment"
$ws.Range("H16").Value = "DDiscount"

$ws.Range("C17").Value = "EQuantity"
$ws.Range("D17").Value = "EPayment"
$ws.Range("E17").Value = "EDiscount"
$ws.Range("F17").Value = "EQuantity"
$ws.Range("G17").Value = "EPayment"
$ws.Range("H17").Value = "EDiscount"

# Anchor cell carries the real array formula, spilling C13:H17.
$ws.Range("C13:H17").FormulaArray = "=C3:C7&D2:I2"

# Columns C:I now share one uniform width instead of per-column bestFit widths.
$ws.Range("C1:I1").EntireColumn.ColumnWidth = 10.296875

# Selection follows the bottom of the new (shorter) spill range.
$ws.Range("I17").Select()
